# [Fonds de solidarite] Add 2021-01-19 data
#
# Refresh nombre_aides (C), nombre_entreprises (D) and montant_total (E)
# for the rows whose figures moved with the 2021-01-19 data update.
# The sheet stores these numeric-looking columns as text, so each cell
# is forced to Text format before the new value is written - this keeps
# the exact textual representation (e.g. trailing ".00") instead of
# letting Excel coerce the assignment into a floating point number.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $r = $ws.Range($range)
    $r.NumberFormat = "@"
    $r.Value = $value
}

# Row 10 - Bourgogne-Franche-Comte / 1 ou 2 salaries
Set-TextValue "C10" "500"
Set-TextValue "D10" "442"
Set-TextValue "E10" "2946256.06"

# Row 11 - Bourgogne-Franche-Comte / 3 a 5 salaries
Set-TextValue "C11" "233"
Set-TextValue "D11" "199"
Set-TextValue "E11" "2173031.14"

# Row 12 - Bourgogne-Franche-Comte / 6 a 9 salaries
Set-TextValue "C12" "73"
Set-TextValue "E12" "1093249.89"

# Row 13 - Bourgogne-Franche-Comte / 10 a 19 salaries
Set-TextValue "C13" "24"
Set-TextValue "E13" "578777.00"

# Row 30 - Corse / 1 ou 2 salaries
Set-TextValue "C30" "180"
Set-TextValue "E30" "788748.58"

# Row 31 - Corse / 3 a 5 salaries
Set-TextValue "C31" "76"
Set-TextValue "D31" "72"
Set-TextValue "E31" "468729.02"

# Row 32 - Corse / 6 a 9 salaries
Set-TextValue "C32" "17"
Set-TextValue "E32" "125000.00"

# Row 74 - La Reunion / 20 a 49 salaries
Set-TextValue "C74" "5"
Set-TextValue "E74" "185000.00"

# Row 93 - Nouvelle-Aquitaine / 1 ou 2 salaries
Set-TextValue "C93" "1116"
Set-TextValue "D93" "1012"
Set-TextValue "E93" "6102866.71"

# Row 95 - Nouvelle-Aquitaine / 6 a 9 salaries
Set-TextValue "C95" "194"
Set-TextValue "E95" "2315043.41"
